$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- String-valued columns, column-major order (matches shared-string allocation order) ---
# Column B
$ws.Range("B210").Value = "NiCrCoFeMn Al0.2"
$ws.Range("B211").Value = "NiCrCoFeMn Al0.6"
$ws.Range("B212").Value = "NiCrCoFeMn Al1.0"
$ws.Range("B213").Value = "NiCrCoFeMn Al1.2"
$ws.Range("B214").Value = "NiCrCoFeMn Al1.6"
$ws.Range("B215").Value = "NiCrCoFeMn Al2.0"
$ws.Range("B216").Value = "NiCrCoFeMn Al0.2"
$ws.Range("B217").Value = "NiCrCoFeMn Al0.6"
$ws.Range("B218").Value = "NiCrCoFeMn Al1.2"
$ws.Range("B219").Value = "NiCrCoFeMn Al1.6"
$ws.Range("B220").Value = "NiCrCoFeMn Al2.0"

# Column C
$ws.Range("C210").Value = "FCC"
$ws.Range("C211").Value = "FCC+BCC+B2"
$ws.Range("C212").Value = "FCC+BCC+B2"
$ws.Range("C213").Value = "FCC+BCC+B2"
$ws.Range("C214").Value = "BCC+FCC+B2"
$ws.Range("C215").Value = "BCC+B2"
$ws.Range("C216").Value = "FCC"
$ws.Range("C217").Value = "FCC+BCC+B2"
$ws.Range("C218").Value = "FCC+BCC+B2"
$ws.Range("C219").Value = "BCC+FCC+B2"
$ws.Range("C220").Value = "BCC+B2"

# Column D
$ws.Range("D210").Value = "SHS"
$ws.Range("D211").Value = "SHS"
$ws.Range("D212").Value = "SHS"
$ws.Range("D213").Value = "SHS"
$ws.Range("D214").Value = "SHS"
$ws.Range("D215").Value = "SHS"
$ws.Range("D216").Value = "SHS"
$ws.Range("D217").Value = "SHS"
$ws.Range("D218").Value = "SHS"
$ws.Range("D219").Value = "SHS"
$ws.Range("D220").Value = "SHS"

# Column E
$ws.Range("E210").Value = "SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"
$ws.Range("E211").Value = "B2 is NiAl; SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"
$ws.Range("E212").Value = "B2 is NiAl; SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"
$ws.Range("E213").Value = "B2 is NiAl; SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"
$ws.Range("E214").Value = "B2 is NiAl; SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"
$ws.Range("E215").Value = "B2 is NiAl; SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"
$ws.Range("E216").Value = "SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"
$ws.Range("E217").Value = "B2 is NiAl; SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"
$ws.Range("E218").Value = "B2 is NiAl; SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"
$ws.Range("E219").Value = "B2 is NiAl; SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"
$ws.Range("E220").Value = "B2 is NiAl; SHS thermite reaction of CoCrFeNiMn oxides with Al into a graphite mold"

# Column F
$ws.Range("F210").Value = "density"
$ws.Range("F211").Value = "density"
$ws.Range("F212").Value = "density"
$ws.Range("F213").Value = "density"
$ws.Range("F214").Value = "density"
$ws.Range("F215").Value = "density"
$ws.Range("F216").Value = "density"
$ws.Range("F217").Value = "density"
$ws.Range("F218").Value = "density"
$ws.Range("F219").Value = "density"
$ws.Range("F220").Value = "density"

# Column G
$ws.Range("G210").Value = "EXP"
$ws.Range("G211").Value = "EXP"
$ws.Range("G212").Value = "EXP"
$ws.Range("G213").Value = "EXP"
$ws.Range("G214").Value = "EXP"
$ws.Range("G215").Value = "EXP"
$ws.Range("G216").Value = "EXP"
$ws.Range("G217").Value = "EXP"
$ws.Range("G218").Value = "EXP"
$ws.Range("G219").Value = "EXP"
$ws.Range("G220").Value = "EXP"

# Column L
$ws.Range("L210").Value = "kg/m^3"
$ws.Range("L211").Value = "kg/m^3"
$ws.Range("L213").Value = "kg/m^3"
$ws.Range("L214").Value = "kg/m^3"
$ws.Range("L215").Value = "kg/m^3"
$ws.Range("L216").Value = "Pa"
$ws.Range("L217").Value = "Pa"
$ws.Range("L218").Value = "Pa"
$ws.Range("L219").Value = "Pa"
$ws.Range("L220").Value = "Pa"

# Column M
$ws.Range("M210").Value = "F2a"
$ws.Range("M211").Value = "F2a"
$ws.Range("M213").Value = "F2a"
$ws.Range("M214").Value = "F2a"
$ws.Range("M215").Value = "F2a"
$ws.Range("M216").Value = "F2b"
$ws.Range("M217").Value = "F2b"
$ws.Range("M218").Value = "F2b"
$ws.Range("M219").Value = "F2b"
$ws.Range("M220").Value = "F2b"

# Column N
$ws.Range("N210").Value = "10.1134/S001250161610002X"
$ws.Range("N211").Value = "10.1134/S001250161610002X"
$ws.Range("N212").Value = "10.1134/S001250161610002X"
$ws.Range("N213").Value = "10.1134/S001250161610002X"
$ws.Range("N214").Value = "10.1134/S001250161610002X"
$ws.Range("N215").Value = "10.1134/S001250161610002X"
$ws.Range("N216").Value = "10.1134/S001250161610002X"
$ws.Range("N217").Value = "10.1134/S001250161610002X"
$ws.Range("N218").Value = "10.1134/S001250161610002X"
$ws.Range("N219").Value = "10.1134/S001250161610002X"
$ws.Range("N220").Value = "10.1134/S001250161610002X"

# --- Numeric columns I, P, Q ---
# Column I
$ws.Range("I210").Value = 298
$ws.Range("I211").Value = 298
$ws.Range("I212").Value = 298
$ws.Range("I213").Value = 298
$ws.Range("I214").Value = 298
$ws.Range("I215").Value = 298
$ws.Range("I216").Value = 298
$ws.Range("I217").Value = 298
$ws.Range("I218").Value = 298
$ws.Range("I219").Value = 298
$ws.Range("I220").Value = 298

# Column P
$ws.Range("P216").Value = 171
$ws.Range("P217").Value = 336
$ws.Range("P218").Value = 382
$ws.Range("P219").Value = 406
$ws.Range("P220").Value = 458

# Column Q
$ws.Range("Q216").Value = 15
$ws.Range("Q217").Value = 25
$ws.Range("Q218").Value = 30
$ws.Range("Q219").Value = 30
$ws.Range("Q220").Value = 35

# --- Plain numeric J column (rows without P/Q-derived formula) ---
$ws.Range("J210").Value = 7700
$ws.Range("J211").Value = 7110
$ws.Range("J213").Value = 6740
$ws.Range("J214").Value = 6610
$ws.Range("J215").Value = 6300

# --- Shared formula for J216:K220 (requires P/Q already populated) ---
$ws.Range("J216:K220").Formula = "=P216*9807000"

# --- Selection / view state to match target ---
$ws.Range("N223").Select()